# Add a new "2022-Q4" quarter sheet to the workbook and record it on the
# "总计" (Total) summary sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert a new worksheet right after "总计" (i.e. before the current
#    2nd tab, "2022-Q2") and name it "2022-Q4". Cloning "2022-Q2" (via
#    Worksheet.Copy) gives the new sheet the same header/style layout
#    every quarterly sheet uses, for free.
# ---------------------------------------------------------------------
$q2Sheet = $wb.Worksheets.Item(2)          # "2022-Q2"
$q2Sheet.Copy($q2Sheet, $null)              # clone placed right before it
$q4Sheet = $wb.Worksheets.Item(2)
$q4Sheet.Name = "2022-Q4"

# Row 2 - 前海开源沪港深裕鑫灵活配置混合C
$q4Sheet.Cells.Item(2,2).Value = "'004317"
$q4Sheet.Cells.Item(2,3).Value = "前海开源沪港深裕鑫灵活配置混合C"
$q4Sheet.Cells.Item(2,4).Value = "'2.88"
$q4Sheet.Cells.Item(2,5).Value = "'90.85"
$q4Sheet.Cells.Item(2,6).Value = "'3.09"
$q4Sheet.Cells.Item(2,7).Value = "'0.0890"
$q4Sheet.Cells.Item(2,8).Value = 4

# Row 3 - 前海开源沪港深裕鑫灵活配置混合A
$q4Sheet.Cells.Item(3,2).Value = "'004316"
$q4Sheet.Cells.Item(3,3).Value = "前海开源沪港深裕鑫灵活配置混合A"
$q4Sheet.Cells.Item(3,4).Value = "'2.30"
$q4Sheet.Cells.Item(3,5).Value = "'90.85"
$q4Sheet.Cells.Item(3,6).Value = "'3.09"
$q4Sheet.Cells.Item(3,7).Value = "'0.0711"
$q4Sheet.Cells.Item(3,8).Value = 4

# ---------------------------------------------------------------------
# 2. Update the "总计" summary sheet: insert a row for 2022-Q4 at the
#    top of the data (row 2), pushing the existing quarters down.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)

# Shift the existing data rows down by one (keeps formatting/styles).
$totalSheet.Range("A2:D6").Copy($totalSheet.Range("A3:D7"))

# Fix up the sequential index column for the shifted rows.
$totalSheet.Cells.Item(3,1).Value = 1
$totalSheet.Cells.Item(4,1).Value = 2
$totalSheet.Cells.Item(5,1).Value = 3
$totalSheet.Cells.Item(6,1).Value = 4
$totalSheet.Cells.Item(7,1).Value = 5

# Write the new 2022-Q4 row.
$totalSheet.Cells.Item(2,1).Value = 0
$totalSheet.Cells.Item(2,2).Value = "2022-Q4"
$totalSheet.Cells.Item(2,3).Value = 2
$totalSheet.Cells.Item(2,4).Value = 0.16

# ---------------------------------------------------------------------
# 3. Restore the originally-active tab ("2020-Q4", the last sheet) since
#    adding a new sheet made it the active one.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Activate()
